$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record (Wins/Losses/Ties), mirroring the
# existing header styling (bold, centered, top-aligned, thin border).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRng = $ws.Range("AD1:AF1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1
$headerRng.Borders.Weight = 2

# Fill every data row (2-52) with the team's record: 73 wins, 89 losses, 0 ties.
$ws.Range("AD2:AD52").Value = 73
$ws.Range("AE2:AE52").Value = 89
$ws.Range("AF2:AF52").Value = 0
